# Preventing kernel tasks from being scheduled on CPU4
# Update the project-plan "Milestones" column (D) with the new weekly
# entries that track the CPU-isolation / Preempt_RT work, and mark the
# two additional weeks (4.3.2024 and 11.3.2024) as completed to match
# the already-highlighted rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Milestones (column D) text for rows 7-15 ---------------------
$ws.Range("D7").Value  = "Trace-cmd und KernelShark funktionieren"
$ws.Range("D8").Value  = "Isolate CPUs of Host"
$ws.Range("D9").Value  = "Decrease latency "
$ws.Range("D10").Value = "Preempt_RT Kernel Patch"
$ws.Range("D11").Value = "Zeitmessungen und Vergleiche zwischen verschiedenen Virtualisierungsmöglichkeiten unter Ubuntu sollen abgeschlossen sein"
$ws.Range("D12").Value = "QEMU/WSL unter Windows lauffähig"
$ws.Range("D13").Value = "Zeitmessungen und Vergleiche zwischen verschiedenen Virualisierungsmöglichkeiten unter Windows sollen abgeschlossen sein"
$ws.Range("D14").Value = "Zeitmessungen abgeschlossen Konklusio und Aufarbeitung/Vergleich PreemptRT gegen Xenomai"
$ws.Range("D15").Value = "Dedizierte Ressourcenzuteilung unter Windows und Messung Verhalten (Core-Sperrung,…)"

# --- Mark weeks 10 and 11 (rows 9-10) as completed, like rows 2-8 --------
# Green fill (matches the existing "done" rows), thin light-grey border.
$greenFill  = 1758337   # RGB(129,212,26)  -> stored as BGR by COM
$borderGrey = 13027014  # RGB(198,198,198) -> stored as BGR by COM

$doneRows = $ws.Range("A9:D10")
$doneRows.Interior.Color = $greenFill
$doneRows.Borders.LineStyle = 1
$doneRows.Borders.Color = $borderGrey

$ws.Range("A9:A10").HorizontalAlignment = -4131
$ws.Range("C9:D10").HorizontalAlignment = -4131
$ws.Range("B9:B10").HorizontalAlignment = -4152
$ws.Range("B9:B10").NumberFormat = "General"

# --- Misc view-state tweaks recorded by the workbook ----------------------
$ws.Range("E11").Select()
$ws.Columns.ColumnWidth = $ws.Columns.ColumnWidth
